$wb = $excel.ActiveWorkbook
$wsCal = $wb.Worksheets.Item("Calificaciones")

$wsCal.Range("J4").Value = 10
$wsCal.Range("J5").Value = 7
$wsCal.Range("X5").Value = 8
$wsCal.Range("J6").Value = 10
$wsCal.Range("J7").Value = 10
$wsCal.Range("X7").Value = 9
$wsCal.Range("J8").Value = 7
$wsCal.Range("X8").Value = 9
$wsCal.Range("J9").Value = 10
$wsCal.Range("J10").Value = 10
$wsCal.Range("J11").Value = 10
$wsCal.Range("X11").Value = 10
$wsCal.Range("J12").Value = 9
$wsCal.Range("X12").Value = 8
$wsCal.Range("J13").Value = 10
$wsCal.Range("J14").Value = 10
$wsCal.Range("J15").Value = 9
$wsCal.Range("J16").Value = 7
$wsCal.Range("X16").Value = 9
$wsCal.Range("J17").Value = 10
$wsCal.Range("X17").Value = 9
$wsCal.Range("J18").Value = 10
$wsCal.Range("J19").Value = 10
$wsCal.Range("J20").Value = 9
$wsCal.Range("J21").Value = 10
$wsCal.Range("J22").Value = 10
$wsCal.Range("J23").Value = 9
$wsCal.Range("X23").Value = 9
$wsCal.Range("J24").Value = 10
$wsCal.Range("X24").Value = 10
$wsCal.Range("J25").Value = 10
$wsCal.Range("J26").Value = 10
$wsCal.Range("J27").Value = 10
$wsCal.Range("X27").Value = 9
$wsCal.Range("J28").Value = 10
$wsCal.Range("J29").Value = 10
$wsCal.Range("X29").Value = 8
$wsCal.Range("J30").Value = 5
$wsCal.Range("X30").Value = 7
$wsCal.Range("J31").Value = 10
$wsCal.Range("X31").Value = 8
$wsCal.Range("J32").Value = 10
$wsCal.Range("X32").Value = 10
$wsCal.Range("J33").Value = 10
$wsCal.Range("J34").Value = 10
$wsCal.Range("X34").Value = 8
$wsCal.Range("J35").Value = 10
$wsCal.Range("J36").Value = 10
$wsCal.Range("X36").Value = 10
$wsCal.Range("J37").Value = 10
$wsCal.Range("J38").Value = 7
$wsCal.Range("X38").Value = 7
$wsCal.Range("J39").Value = 5
$wsCal.Range("X39").Value = 7
$wsCal.Range("J40").Value = 10
$wsCal.Range("X40").Value = 10
$wsCal.Range("J41").Value = 10
$wsCal.Range("J42").Value = 10
$wsCal.Range("X42").Value = 8
$wsCal.Range("J43").Value = 10

$wsTot = $wb.Worksheets.Item("Totales")
$wsTot.Range("H6").Value = 9.300000000000001
